# Renamed few transcripts. Updated the DataSheet
# Replace every "RBD" speaker tag in column D with "T".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,8,10,16,24,25,26,27,28,29,30,31,34,36,37,40,41,42,44)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 4)  # Column D
    if ($cell.Value2 -eq "RBD") {
        $cell.Value = "T"
    }
}
